# feat: add setWidth, setHeight
#
# Resize the embedded "Score Transition" chart on Sheet1 by growing its
# width and height (keeping its top-left anchor where it is). This moves
# the drawing's bottom-right (two-cell) anchor from col 10 / row 14 out to
# col 13 / row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$chartObj = $ws.ChartObjects(1)

$chartObj.Width  = 584.375
$chartObj.Height = 360
